$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("A1")
    if ($cell.Value2 -eq "Month/Year: FEBRUARY 2021") {
        $cell.Value = "Month/Year: APRIL 2021"
        # Setting the value recalculates the row's autofit height because
        # A1 uses a much larger font; restore the natural (default) row
        # height so no explicit row height survives in the saved file.
        $ws.Rows.Item(1).AutoFit()
    }
}
